# Apply crypto price/volume updates scraped on Mon May 15 11:59:35 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $ws.Range("B2").Style
}

Set-TextCell "D2" "27.675.16"
$ws.Range("E2").Value = "  -0.32%  "
Set-TextCell "D3" "1.848.53"
$ws.Range("E3").Value = "  -0.82%  "
Set-TextCell "D5" "319.72"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("E6").Value = "  -2.30%  "
Set-TextCell "D7" "0.4307"
Set-TextCell "D8" "0.3743"
$ws.Range("E8").Value = "  -1.41%  "
Set-TextCell "D9" "0.07347"
$ws.Range("E9").Value = "  -1.51%  "
Set-TextCell "D10" "0.8799"
$ws.Range("E10").Value = "  -0.54%  "
Set-TextCell "D11" "21.56"
$ws.Range("E11").Value = "  -0.55%  "
Set-TextCell "D12" "1.858.16"
$ws.Range("E12").Value = "  -0.52%  "
Set-TextCell "D13" "6.725"
$ws.Range("E13").Value = "  -0.58%  "
Set-TextCell "D14" "5.453"
$ws.Range("E14").Value = "  -1.88%  "
Set-TextCell "D15" "0.07129"
$ws.Range("E15").Value = "  -1.39%  "
Set-TextCell "D16" "87.81"
$ws.Range("E16").Value = "  +4.83%  "
$ws.Range("E17").Value = "  -2.54%  "
Set-TextCell "D18" "0.000008995"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("E20").Value = "  -0.47%  "
Set-TextCell "D21" "27.676.00"
$ws.Range("E21").Value = "  -0.35%  "
Set-TextCell "D22" "5.242"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("E23").Value = "  -1.81%  "
Set-TextCell "D24" "2.080.20"
$ws.Range("E24").Value = "  -0.44%  "
Set-TextCell "D25" "2.011"
$ws.Range("E25").Value = "  +0.25%  "
Set-TextCell "D26" "155.90"
$ws.Range("E26").Value = "  -1.81%  "
Set-TextCell "D27" "18.62"
$ws.Range("E27").Value = "  -1.28%  "
Set-TextCell "D28" "2.122"
$ws.Range("E28").Value = "  +7.05%  "
Set-TextCell "D29" "5.387"
$ws.Range("E29").Value = "  +1.18%  "
Set-TextCell "D30" "120.53"
$ws.Range("E30").Value = "  +2.08%  "
Set-TextCell "D31" "0.08929"
$ws.Range("E31").Value = "  -1.51%  "
Set-TextCell "D32" "1.225"
$ws.Range("E32").Value = "  +0.87%  "
Set-TextCell "D33" "0.7778"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -0.19%  "
Set-TextCell "D35" "2.925"
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("E36").Value = "  -2.52%  "
Set-TextCell "D37" "1.140"
$ws.Range("E37").Value = "  -0.94%  "
Set-TextCell "D38" "0.05341"
$ws.Range("E38").Value = "  -0.12%  "
Set-TextCell "D39" "0.01973"
$ws.Range("E39").Value = "  -1.06%  "
Set-TextCell "D40" "7.192"
$ws.Range("E40").Value = "  +4.41%  "
Set-TextCell "D41" "2.879"
$ws.Range("E41").Value = "  +0.68%  "
Set-TextCell "D42" "0.5157"
$ws.Range("E42").Value = "  -0.84%  "
Set-TextCell "D43" "0.1682"
$ws.Range("E43").Value = "  -0.67%  "
Set-TextCell "D44" "8.870"
$ws.Range("E44").Value = "  +2.33%  "
Set-TextCell "B45" "Quant"
Set-TextCell "C45" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D45" "108.87"
$ws.Range("E45").Value = "  -0.83%  "
Set-TextCell "B46" "EnergySwap"
Set-TextCell "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D46" "10.63"
$ws.Range("E46").Value = "  -0.98%  "
Set-TextCell "D47" "0.4732"
$ws.Range("E47").Value = "  +0.50%  "
Set-TextCell "D48" "0.06505"
$ws.Range("E48").Value = "  +0.43%  "
Set-TextCell "D49" "1.700"
$ws.Range("E49").Value = "  -1.16%  "
Set-TextCell "D50" "1.013"
Set-TextCell "D51" "1.875"
$ws.Range("E51").Value = "  -2.23%  "
